# Regenerate the save_data "K" column (strikeouts, formerly "Strike#") for
# each outing row (rows 2-69 on Sheet1) with freshly-computed values, as part
# of regenerating std/mean and the s_vals used downstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, keyed by row number (row 2 = most recent outing ... row 69 =
# oldest). These are the recomputed "K" counts for each game.
$kValues = [ordered]@{
    2  = 0
    3  = 2
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 3
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    41 = 1
    42 = 3
    43 = 0
    44 = 2
    45 = 1
    46 = 2
    47 = 1
    48 = 3
    49 = 2
    50 = 1
    51 = 1
    52 = 2
    53 = 2
    54 = 0
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    59 = 2
    60 = 2
    61 = 2
    62 = 2
    63 = 2
    64 = 2
    65 = 3
    66 = 2
    67 = 0
    68 = 1
    69 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

Write-Output "Updated K column (G2:G69) with regenerated s_vals"
